# Switched to EPPlus instead of Excel.Interop
$wb = $excel.ActiveWorkbook

# Rename the "Process" sheet to "burp"
$wsProcess = $wb.Worksheets.Item("Process")
$wsProcess.Name = "burp"

# Update the Input sheet: D4 changes from 8 to 9, selection moves to E4
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("D4").Value = 9
$wsInput.Range("E4").Select() | Out-Null

# Update selection on the renamed "burp" sheet
$wsBurp = $wb.Worksheets.Item("burp")
$wsBurp.Range("B5").Select() | Out-Null

# Update the Output sheet formula for C5 to reference the renamed sheet
# and add the explicit FALSE range-lookup argument
$wsOutput = $wb.Worksheets.Item("Output")
$wsOutput.Range("C5").Formula = "=((VLOOKUP(Input!D5,burp!A2:B5,2,FALSE))*C3)/1000"

# Restore Output as the active sheet (it was active before these edits)
$wsOutput.Activate() | Out-Null

$wb.Application.Calculate()
